$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Oklo Inc. / OKLO)
$ws.Range("D2").Value = 111.65
$ws.Range("E2").Value = 57.4
$ws.Range("F2").Value = 25.85
$ws.Range("H2").Value = 73
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 83
$ws.Range("K2").Value = 59.5
$ws.Range("M2").Value = "⛔ 관망하십시오."
$ws.Range("N2").Value = 54.85170003294819
$ws.Range("O2").Value = "⚪ 중립 구간"

# Row 3 (NuScale Power Corporation / SMR)
$ws.Range("D3").Value = 22.85
$ws.Range("E3").Value = 49.1
$ws.Range("F3").Value = 20.01
$ws.Range("K3").Value = 53.9
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 54.85170003294819
$ws.Range("O3").Value = "⚪ 중립 구간"
